$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Last
$r = $newp.Range

$wordml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>باگ</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>کد</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>مج</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>د</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>پس‌ترت</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>ب</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>خانه</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>هزار</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>اتاق</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>مساحت</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>محصور</w:t></w:r></w:p>
'@

$xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $wordml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$r.InsertXML($xml)

Write-Host "Count after insert: $($d.Paragraphs.Count)"

# Merge away the spare trailing empty paragraph created by InsertXML,
# by deleting the paragraph mark of the second-to-last paragraph.
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$markRange = $d.Range($secondLast.Range.End - 1, $secondLast.Range.End)
$markRange.Delete()

Write-Host "Count after merge: $($d.Paragraphs.Count)"
